# sdf_deja_traite.xlsx - "Encore des chgmts de titres"
#
# The "Year of Treatment" column (B) is removed entirely, shifting the
# accommodation-status columns (old C:H) one place to the left (new B:G).
# The new header row (B1:G1) additionally gets a ".deja.deja.deja" suffix
# appended to its title text (the "Country" header in A1 is left alone).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Year of Treatment" column; everything to its right shifts left.
$ws.Columns("B").Delete()

# Append the ".deja.deja.deja" suffix to the (now shifted) header titles.
$suffix = ".deja.deja.deja"
foreach ($col in @("B", "C", "D", "E", "F", "G")) {
    $cell = $ws.Range($col + "1")
    $cell.Value = $cell.Value() + $suffix
}
